$wb = $excel.ActiveWorkbook


$headerIND = @("Sector","Year","COL","OIL","GAS","BMS","ELE")
$headerELE = @("Year","COL","COLX","OIL","OILX","GAS","GASX","NUC","BMS","BMSX","HYD","GEO","WIN","PV")
$headerLOSS = @("Year","ELE")

# LOSS: 11 rows
$LOSS_data = @(
@(2010,0.905095428809249),
@(2011,0.906341796876974),
@(2012,0.90931390033079),
@(2013,0.905790118440597),
@(2014,0.910844163083025),
@(2015,0.911547003404931),
@(2016,0.897322110943927),
@(2017,0.899794083235142),
@(2018,0.8983272718829),
@(2019,0.898813117102369),
@(2020,0.89526829442089)
)

# sIND: 11 rows
$sIND_data = @(
@("Industry",2010,0.235317279731797,0.2542029283506,0.118208827636883,0.0364559087844255,0.355815055496294),
@("Industry",2011,0.233781164697422,0.265632827450878,0.123996041168318,0.0393771825968263,0.337212784086556),
@("Industry",2012,0.238183139592608,0.254729802181916,0.126204574849316,0.0385824030204225,0.342300080355736),
@("Industry",2013,0.238991514104705,0.25628028219994,0.118880058675235,0.0418285950752597,0.34401954994486),
@("Industry",2014,0.239516304821693,0.252808385648269,0.119951336210655,0.0394015934257564,0.348322379893626),
@("Industry",2015,0.244278082479678,0.240704426350651,0.123588729268783,0.040869790218956,0.350558971681932),
@("Industry",2016,0.242373154334521,0.228317868362445,0.130553448532646,0.0424868181183189,0.356268710652069),
@("Industry",2017,0.24922368028528,0.217978539482308,0.130354331207477,0.0435277511261819,0.358915697898754),
@("Industry",2018,0.242436419524306,0.217749202182963,0.136550016297178,0.0452056222686474,0.358058739726906),
@("Industry",2019,0.246172595109221,0.21124427442343,0.134886409911155,0.0459790921710079,0.361717628385187),
@("Industry",2020,0.2378380422764,0.218958939914931,0.134617175435078,0.0445472571891143,0.364038585184477)
)

# sTRA: 11 rows
$sTRA_data = @(
@("Transport",2010,0.0000127563610961359,0.976242188627097,0.00127593983249683,0.00252089993090304,0.0199482152484074),
@("Transport",2011,0.0000130761890381062,0.97628601984156,0.00128582525541378,0.00251280766015607,0.0199022710538319),
@("Transport",2012,0.0000122700235773222,0.976158085725306,0.00122448542981893,0.00258677266289086,0.0200183861584066),
@("Transport",2013,0.0000120951891385202,0.97548941750524,0.00113121847890265,0.00306135602984965,0.0203059127968693),
@("Transport",2014,0.0000124358001815627,0.974323345204068,0.00109827751077169,0.00382662661902664,0.0207393148659519),
@("Transport",2015,0.0000121347680783447,0.97356260392445,0.000978980613888078,0.00472436038293392,0.0207219203106501),
@("Transport",2016,0.0000132625027271021,0.972840052227736,0.000840511110330097,0.00536667172852187,0.0209395024306852),
@("Transport",2017,0.0000133731630706788,0.972273086676485,0.000719476173202521,0.00587984547310072,0.0211142185141413),
@("Transport",2018,0.0000121864816036597,0.972229378018905,0.000593752464800533,0.00607733067306953,0.0210873523616217),
@("Transport",2019,0.0000124402436352604,0.972049191487832,0.000483441690159145,0.00592846721684796,0.0215264593615252),
@("Transport",2020,0.0000137654722952663,0.968745495153598,0.000406463806940781,0.00701465525712946,0.0238196203100367)
)

# sCOM: 11 rows
$sCOM_data = @(
@("Commercial and public services",2010,0.00465182109398443,0.248780292784875,0.163919367208958,0.025323313743574,0.557325205168608),
@("Commercial and public services",2011,0.00443937837618143,0.247661613325155,0.170166010212602,0.0269752739741212,0.55075772411194),
@("Commercial and public services",2012,0.00370231464474871,0.241225779086881,0.168257442302387,0.028215585444752,0.558598878521231),
@("Commercial and public services",2013,0.00377392792263585,0.267811848087296,0.166527881251453,0.0293294768802379,0.532556865858377),
@("Commercial and public services",2014,0.0036354345163539,0.255314983044822,0.167067980866524,0.0296629979050195,0.54431860366728),
@("Commercial and public services",2015,0.00384221225809996,0.235266941064898,0.186677432377876,0.030863828882524,0.543349585416602),
@("Commercial and public services",2016,0.00344430509299673,0.22954751401036,0.172581337311204,0.0318521241103312,0.562574719475109),
@("Commercial and public services",2017,0.00352644182596447,0.216322072075319,0.186930731705442,0.0336621891716418,0.559558565221633),
@("Commercial and public services",2018,0.00227920610692883,0.240628984112254,0.176965719300653,0.0340575770647568,0.546068513415407),
@("Commercial and public services",2019,0.00254169586888653,0.215081088233997,0.188233506667908,0.0355763905566147,0.558567318672594),
@("Commercial and public services",2020,0.00269743384889377,0.235100657037878,0.155260048019108,0.0375768243307909,0.569365036763329)
)

# sRES: 11 rows
$sRES_data = @(
@("Residential",2010,0,0.295577388456578,0.185649697865496,0.000194660599825144,0.518578253078101),
@("Residential",2011,0,0.294781866130855,0.192891623759889,0.000179329398467914,0.512147180710787),
@("Residential",2012,0,0.292714418499285,0.190718860799803,0.000168908764924207,0.516397811935988),
@("Residential",2013,0,0.288747412097999,0.191014059935836,0.000173333992682247,0.520065193973483),
@("Residential",2014,0,0.282634731811053,0.199357563841297,0.000180578946725493,0.517827125400925),
@("Residential",2015,0,0.276882966712301,0.196515236635559,0.000185446790406947,0.526416349861733),
@("Residential",2016,0,0.275029440633412,0.200188136813701,0.00018519547382262,0.524597227079064),
@("Residential",2017,0,0.282907320193399,0.201477249685925,0.000177679670728215,0.515437750449948),
@("Residential",2018,0,0.261349883050479,0.204412917935886,0.000192649405279954,0.534044549608355),
@("Residential",2019,0,0.272207345338359,0.208935466078323,0.000194336241708559,0.518662852341609),
@("Residential",2020,0,0.267849508803155,0.212616037620618,0.00018539845399501,0.519349055122232)
)

# sELE: 11 rows
$sELE_data = @(
@(2010,0.277508981658784,0,0.0794383653074846,0,0.290670709251437,0,0.252133772699081,0.0180280502430537,0,0.0733053816472367,0.00230238382452896,0.00351305981736638,0.00309929555102816),
@(2011,0.279514031567449,0,0.138739425883535,0,0.378150073272013,0,0.093966393678015,0.0199593518081612,0,0.0784135200946303,0.00247102592822759,0.00431783155470561,0.00446834621326357),
@(2012,0.309036885418116,0,0.16871896917671,0,0.401733664142025,0,0.0148835897660963,0.0210362952663667,0,0.0714615873774758,0.00243624353471017,0.00451764899230655,0.00617511632619331),
@(2013,0.32857950846392,0,0.137454759567475,0,0.410167848157459,0,0.0086399126257143,0.0222308490295306,0,0.0737313407885046,0.00241746668437432,0.00481728762652694,0.0119610270564952),
@(2014,0.334326080980433,0,0.101734914500362,0,0.431595057885971,0,0,0.0236005064249575,0,0.0794565919426686,0.00249217390808064,0.00496246995360942,0.0218322044039186),
@(2015,0.341673503664889,0,0.0884855919937461,0,0.410509175767614,0,0.00913029512612327,0.0243335861732676,0,0.0842866430854728,0.00251066184722792,0.00539864859635136,0.0336718937453076),
@(2016,0.330675581183728,0,0.0705485684728781,0,0.420809586298719,0,0.0172990162741022,0.0323968620389084,0,0.0761357867412523,0.00239561681625303,0.00590618684087012,0.0438327953332886),
@(2017,0.331369893457585,0,0.0577521051101989,0,0.403932324957214,0,0.0312580134749847,0.0353048777013127,0,0.0795838976402631,0.00233352391553346,0.00616384624005379,0.052301517502854),
@(2018,0.319154245328742,0,0.0437033390053524,0,0.389648517165296,0,0.0623144717104449,0.0377251652419438,0,0.0777076213317824,0.00242236483847222,0.00717975885761121,0.0601445165203554),
@(2019,0.322438841124481,0,0.0343351533002472,0,0.381051176338918,0,0.0628909173560157,0.0416705698826866,0,0.0788791475573428,0.00280735730746134,0.00751093804388234,0.0684158990889647),
@(2020,0.313426983789406,0,0.0321379578663662,0,0.397941195147002,0,0.0390594368272674,0.0462308832823151,0,0.0794322109322477,0.00301573686486334,0.00904116299392519,0.0797144322966066)
)


# --- Step 1: repurpose the existing "sIND" sheet into "LOSS" ---
$lossSheet = $wb.Worksheets.Item("sIND")
$lossSheet.Name = "LOSS"
$lossSheet.Cells.Clear()

for ($c = 0; $c -lt $headerLOSS.Length; $c++) {
    $lossSheet.Cells.Item(1, $c + 1).Value = $headerLOSS[$c]
}
for ($i = 0; $i -lt $LOSS_data.Length; $i++) {
    $row = $LOSS_data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $lossSheet.Cells.Item($i + 2, $c + 1).Value = $row[$c]
    }
}

# --- Step 2: add the new "sIND" sheet (industry energy-mix shares) after LOSS ---
$sINDSheet = $wb.Worksheets.Add($null, $lossSheet)
$sINDSheet.Name = "sIND"
for ($c = 0; $c -lt $headerIND.Length; $c++) {
    $sINDSheet.Cells.Item(1, $c + 1).Value = $headerIND[$c]
}
for ($i = 0; $i -lt $sIND_data.Length; $i++) {
    $row = $sIND_data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sINDSheet.Cells.Item($i + 2, $c + 1).Value = $row[$c]
    }
}

# --- Step 3: add "sTRA" sheet (transport energy-mix shares) ---
$sTRASheet = $wb.Worksheets.Add($null, $sINDSheet)
$sTRASheet.Name = "sTRA"
for ($c = 0; $c -lt $headerIND.Length; $c++) {
    $sTRASheet.Cells.Item(1, $c + 1).Value = $headerIND[$c]
}
for ($i = 0; $i -lt $sTRA_data.Length; $i++) {
    $row = $sTRA_data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sTRASheet.Cells.Item($i + 2, $c + 1).Value = $row[$c]
    }
}

# --- Step 4: add "sCOM" sheet (commercial energy-mix shares) ---
$sCOMSheet = $wb.Worksheets.Add($null, $sTRASheet)
$sCOMSheet.Name = "sCOM"
for ($c = 0; $c -lt $headerIND.Length; $c++) {
    $sCOMSheet.Cells.Item(1, $c + 1).Value = $headerIND[$c]
}
for ($i = 0; $i -lt $sCOM_data.Length; $i++) {
    $row = $sCOM_data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sCOMSheet.Cells.Item($i + 2, $c + 1).Value = $row[$c]
    }
}

# --- Step 5: add "sRES" sheet (residential energy-mix shares) ---
$sRESSheet = $wb.Worksheets.Add($null, $sCOMSheet)
$sRESSheet.Name = "sRES"
for ($c = 0; $c -lt $headerIND.Length; $c++) {
    $sRESSheet.Cells.Item(1, $c + 1).Value = $headerIND[$c]
}
for ($i = 0; $i -lt $sRES_data.Length; $i++) {
    $row = $sRES_data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sRESSheet.Cells.Item($i + 2, $c + 1).Value = $row[$c]
    }
}

# --- Step 6: add "sELE" sheet (electricity generation mix by source) ---
$sELESheet = $wb.Worksheets.Add($null, $sRESSheet)
$sELESheet.Name = "sELE"
for ($c = 0; $c -lt $headerELE.Length; $c++) {
    $sELESheet.Cells.Item(1, $c + 1).Value = $headerELE[$c]
}
for ($i = 0; $i -lt $sELE_data.Length; $i++) {
    $row = $sELE_data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sELESheet.Cells.Item($i + 2, $c + 1).Value = $row[$c]
    }
}

Write-Output "Edit complete"
